$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old content (columns A:E, rows 1:6) since the layout changes drastically
$ws.Range("A1:E6").Clear()

# New headers (written in this specific order so the shared-strings table
# builds up in the same sequence the original author typed them: surviving
# strings keep their old relative order, brand-new ones are appended in the
# order they are first used)
$ws.Range("B1").Value = "Cust_name"
$ws.Range("K1").Value = "Total"
$ws.Range("A1").Value = "Sr.No"
$ws.Range("C1").Value = "Cow"
$ws.Range("F1").Value = "Buffalo"
$ws.Range("I1").Value = "Other"
$ws.Range("J1").Value = "Previous_pending"
$ws.Range("E1").Value = "CM_total"
$ws.Range("D1").Value = "C_rate"
$ws.Range("G1").Value = "B_rate"
$ws.Range("H1").Value = "BM_total"

# New data row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Wagvilasinee Kulkarni"
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 50
$ws.Range("E2").Formula = "=D2*C2"
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = 75
$ws.Range("H2").Formula = "=G2*F2"
$ws.Range("I2").Value = 250
$ws.Range("J2").Value = 500
$ws.Range("K2").Formula = "=E2+H2+I2+J2"

# Column widths (closest values the host's quantized ColumnWidth storage
# can reach to the target 12.140625 / 15.7109375 / 17.28515625 "best fit"
# widths recorded in the authored workbook)
$ws.Range("E1").ColumnWidth = 11.333333333333334
$ws.Range("H1").ColumnWidth = 14.833333333333334
$ws.Range("J1").ColumnWidth = 16.5

# Page setup - portrait orientation
$ws.PageSetup.Orientation = 1

# Selection
$ws.Range("H6").Select()
